# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "conversion" note (A1) with the new day's rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$nuevoTexto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.44 = 13233.54 pesos`n✅ 13233.54 pesos = 3.43 = 957.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $nuevoTexto

# --- tasas: update the rate cells N10/O10 and N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 290.7
$ws2.Range("O10").Value = 3846.99

$ws2.Range("N12").Value = 3860.99
$ws2.Range("O12").Value = 279.26
